$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.187.36'
$ws.Range("E2").Value = '  +2.21%  '

$ws.Range("D3").Value = '2.351.80'
$ws.Range("E3").Value = '  +6.46%  '

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.03'
$ws.Range("E5").Value = '  +5.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.07'
$ws.Range("E6").Value = '  +1.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.640'
$ws.Range("E7").Value = '  +3.25%  '

$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.636'
$ws.Range("E9").Value = '  +7.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.03'
$ws.Range("E10").Value = '  -2.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0941'
$ws.Range("E11").Value = '  +3.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.97'
$ws.Range("E12").Value = '  +2.62%  '

$ws.Range("E13").Value = '  +11.78%  '

$ws.Range("E14").Value = '  +2.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.34'
$ws.Range("E15").Value = '  +10.07%  '

$ws.Range("D16").Value = '2.705.76'
$ws.Range("E16").Value = '  +6.59%  '

$ws.Range("D17").Value = '2.339.17'
$ws.Range("E17").Value = '  +4.95%  '

$ws.Range("D18").Value = '43.157.03'
$ws.Range("E18").Value = '  +2.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000109'
$ws.Range("E19").Value = '  +4.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'
$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.55'
$ws.Range("E21").Value = '  +4.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.43'
$ws.Range("E22").Value = '  -0.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.55'
$ws.Range("E23").Value = '  +12.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '253.13'
$ws.Range("E24").Value = '  +11.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.08'
$ws.Range("E25").Value = '  +1.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.04'
$ws.Range("E26").Value = '  +4.37%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.99'
$ws.Range("E28").Value = '  +1.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.21'
$ws.Range("E29").Value = '  +3.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.24'
$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.61'
$ws.Range("E31").Value = '  +8.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.64'
$ws.Range("E32").Value = '  +0.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.16'
$ws.Range("E33").Value = '  -0.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0930'
$ws.Range("E34").Value = '  +7.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.95'
$ws.Range("E35").Value = '  +8.20%  '

$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.133'
$ws.Range("E36").Value = '  +6.06%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.03'
$ws.Range("E37").Value = '  +1.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0379'
$ws.Range("E38").Value = '  +5.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.11'
$ws.Range("E39").Value = '  -4.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +1.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.71'
$ws.Range("E41").Value = '  +10.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.56'
$ws.Range("E42").Value = '  +4.00%  '

$ws.Range("E43").Value = '  +15.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.233'
$ws.Range("E44").Value = '  +1.91%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.64'
$ws.Range("E46").Value = '  -0.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.64'
$ws.Range("E47").Value = '  +4.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.42'
$ws.Range("E48").Value = '  +12.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.35'
$ws.Range("E49").Value = '  +7.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.30'
$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.101'
$ws.Range("E51").Value = '  +3.52%  '
